# Apply the "1bn instruction runs" data update to the exec_time worksheet.
# The chart sheets (bc_exec_time, umc_exec_time, bc_coverage, umc_coverage)
# pull their series values from this sheet via formula references, so
# updating the backing cells is the substantive edit; Excel keeps their
# cached numCache values in sync when the workbook is recalculated/saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exec_time")

# --- BC execution-time overhead-difference table (rows 4:6) ---
# Row 4: replace volatile RAND() formulas with literal values
$ws.Range("B4").Value = 0.04
$ws.Range("C4").Value = 0.0309
$ws.Range("D4").Value = -0.0169
$ws.Range("E4").Value = 0.09
$ws.Range("F4").Value = -0.3501
$ws.Range("G4").Value = -0.5194
$ws.Range("H4").Value = 0.146
$ws.Range("I4").Value = 0.1275

# Row 5
$ws.Range("B5").Value = 15.14
$ws.Range("C5").Value = 14.64
$ws.Range("D5").Value = 13.64
$ws.Range("E5").Value = 12.64
$ws.Range("F5").Value = 11.64
$ws.Range("G5").Value = 10.65
$ws.Range("H5").Value = 9.65
$ws.Range("I5").Value = 8.65

# Row 6: cleared entirely (no formula, no value)
$ws.Range("B6:I6").ClearContents()

# --- UMC execution-time overhead-difference table (rows 10:12) ---
# Row 10
$ws.Range("B10").Value = 0.0115
$ws.Range("C10").Value = 0.0057
$ws.Range("D10").Value = -0.0093
$ws.Range("E10").Value = -0.0256
$ws.Range("F10").Value = -0.4356
$ws.Range("G10").Value = -0.6236
$ws.Range("H10").Value = -1.1329
$ws.Range("I10").Value = -2.0856

# Row 11
$ws.Range("B11").Value = 0.0119
$ws.Range("C11").Value = 0.0065
$ws.Range("D11").Value = -0.0025
$ws.Range("E11").Value = -0.0093
$ws.Range("F11").Value = -0.5835
$ws.Range("G11").Value = -1.5835
$ws.Range("H11").Value = -2.5835
$ws.Range("I11").Value = -3.5835

# Row 12: cleared entirely
$ws.Range("B12:I12").ClearContents()

# --- BC coverage table (rows 17:19) ---
# Row 17
$ws.Range("B17").Value = 0.88
$ws.Range("C17").Value = 0.89
$ws.Range("D17").Value = 0.91
$ws.Range("E17").Value = 0.92
$ws.Range("F17").Value = 0.93
$ws.Range("G17").Value = 0.94
$ws.Range("H17").Value = 0.96
$ws.Range("I17").Value = 0.97

# Row 18
$ws.Range("B18").Value = 0.94
$ws.Range("C18").Value = 0.94
$ws.Range("D18").Value = 0.96
$ws.Range("E18").Value = 0.96
$ws.Range("F18").Value = 0.97
$ws.Range("G18").Value = 0.97
$ws.Range("H18").Value = 0.98
$ws.Range("I18").Value = 1

# Row 19: cleared entirely
$ws.Range("B19:I19").ClearContents()

# --- UMC coverage table (rows 23:25) ---
# Row 23
$ws.Range("B23").Value = 0.49
$ws.Range("C23").Value = 0.63
$ws.Range("D23").Value = 0.82
$ws.Range("E23").Value = 0.95
$ws.Range("F23").Value = 0.98
$ws.Range("G23").Value = 0.99
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1

# Row 24
$ws.Range("B24").Value = 0.52
$ws.Range("C24").Value = 0.67
$ws.Range("D24").Value = 0.88
$ws.Range("E24").Value = 0.99
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1

# Row 25: cleared entirely
$ws.Range("B25:I25").ClearContents()

# --- Sheet view: exec_time becomes the active/selected sheet, scrolled
#     back to the top, with B17:I25 selected (active cell B17) ---
$ws.Activate()
$ws.Range("B17:I25").Select()
